# Update "想去人数" (number of people interested) counts in column F
# for the 展览 (Exhibitions) and 全部类型 (All Types) sheets, matching the
# freshly generated site output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 26
$wsExhibit.Range("F6").Value = 327
$wsExhibit.Range("F9").Value = 109
$wsExhibit.Range("F18").Value = 1686
$wsExhibit.Range("F21").Value = 258
$wsExhibit.Range("F22").Value = 2874
$wsExhibit.Range("F26").Value = 932
$wsExhibit.Range("F29").Value = 2852
$wsExhibit.Range("F35").Value = 1904
$wsExhibit.Range("F37").Value = 1913
$wsExhibit.Range("F39").Value = 37
$wsExhibit.Range("F40").Value = 50

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 327
$wsAll.Range("F10").Value = 109
$wsAll.Range("F19").Value = 1686
$wsAll.Range("F22").Value = 258
$wsAll.Range("F23").Value = 2874
$wsAll.Range("F28").Value = 2852
$wsAll.Range("F36").Value = 1904
$wsAll.Range("F39").Value = 1913
